$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.01"
$ws.Range("D3").Value = "'23.99"
$ws.Range("D4").Value = "'5.401"
$ws.Range("D5").Value = "'0.05897"
$ws.Range("D6").Value = "'3.398"
$ws.Range("D7").Value = "'6.503"
$ws.Range("D8").Value = "'0.8115"
$ws.Range("D9").Value = "'0.9248"
$ws.Range("D10").Value = "'0.1419"
$ws.Range("D11").Value = "'0.07404"
$ws.Range("D12").Value = "'0.03069"
$ws.Range("D13").Value = "'0.03084"
$ws.Range("D14").Value = "'0.09337"
$ws.Range("D15").Value = "'3.865"
$ws.Range("D16").Value = "'0.001552"
$ws.Range("D17").Value = "'0.04720"
$ws.Range("D18").Value = "'0.0005978"
$ws.Range("D19").Value = "'0.005936"
$ws.Range("D20").Value = "'0.001251"
$ws.Range("D21").Value = "'0.004728"
$ws.Range("D22").Value = "'0.00008819"
$ws.Range("D23").Value = "'3.560"
$ws.Range("D27").Value = "'0.0002656"
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"
$ws.Range("D41").Value = "'0.006374"
$ws.Range("D43").Value = "'0.002746"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.008556"
$ws.Range("D45").Value = "'0.00005225"
$ws.Range("D47").Value = "'0.6718"
$ws.Range("D48").Value = "'0.001961"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D50").Value = "'0.0002002"
